$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "trainingimages/24_takopa"
$ws.Range("B2").Value = "pngimages/24_banana.png"
$ws.Range("C2").Value = "trainingimages/15_kopota"
$ws.Range("D2").Value = "pngimages/15_barrel.png"

# Row 3
$ws.Range("A3").Value = "trainingimages/20_tatito"
$ws.Range("B3").Value = "pngimages/20_pizza.png"
$ws.Range("C3").Value = "trainingimages/13_kopopi"
$ws.Range("D3").Value = "pngimages/13_toast.png"
$ws.Range("E3").Value = -0.5
$ws.Range("F3").Value = 0.5

# Row 4
$ws.Range("A4").Value = "trainingimages/02_pitito"
$ws.Range("B4").Value = "pngimages/02_pallet.png"
$ws.Range("C4").Value = "trainingimages/03_kikita"
$ws.Range("D4").Value = "pngimages/03_box.png"
